$d = $word.ActiveDocument

# Update the ID placeholder text in the first paragraph's first run.
$null = $d.Content.Find.Execute("**ID__AFFARS_mp_5301_602_2_d_topic_4__ID**", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_MP5301_602_2_3__ID**", 2)

$p = $d.Paragraphs(1)

# Remove the now-orphaned trailing space run at the end of the paragraph
# (the run just before the paragraph mark).
$spaceRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

# Add the paragraph border (matching the other body paragraphs) and
# widen the left indent from 120 to 225 twips (6pt -> 11.25pt).
$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromRight = 5
$p.LeftIndent = 11.25
